$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@(2, 14, 'Aristocrat Pub & Restaurant', 2, 4.5, 1182)
  ,@(3, 23, 'Axum Ethiopian Restaurant', 2, 4.7, 413)
  ,@(4, 44, 'BRU Burger Bar', 2, 4.6, 4150)
  ,@(5, 53, 'Big Hoffa''s Smokehouse', 1, 4.7, 2284)
  ,@(6, 10, 'Bluebeard', 3, 4.7, 1366)
  ,@(7, 19, 'Bonefish Grill', 2, 4.5, 1081)
  ,@(8, 29, 'Bosphorus Istanbul Cafe', 2, 4.5, 1271)
  ,@(9, 5, 'Burritos & Beer Restaurant, LLC', 1, 4.7, 329)
  ,@(10, 57, 'Carrabba''s Italian Grill', 2, 4.4, 1120)
  ,@(11, 27, 'Charleston''s Restaurant', 2, 4.5, 1058)
  ,@(12, 25, 'Chuy''s', 2, 4.4, 2272)
  ,@(13, 4, 'Cooper''s Hawk Winery & Restaurant', 2, 4.6, 1498)
  ,@(14, 41, 'Courses Restaurant', $null, 4.6, 38)
  ,@(15, 28, 'Cracker Barrel Old Country Store', 2, 4.4, 2750)
  ,@(16, 58, 'Culver''s', 1, 4.5, 84)
  ,@(17, 34, 'Fire by the Monon', 2, 4.6, 906)
  ,@(18, 36, 'First Watch', 2, 4.6, 396)
  ,@(19, 43, 'Flatwater', 2, 4.6, 874)
  ,@(20, 18, 'Greek Islands', 2, 4.6, 866)
  ,@(21, 48, 'His Place Eatery - Chicken & Waffles, Ribs and Soul Food', 1, 4.5, 2164)
  ,@(22, 22, 'Houlihan''s', 2, 4.3, 836)
  ,@(23, 52, 'IHOP', 1, 4.1, 2047)
  ,@(24, 49, 'Iaria''s Italian Restaurant', 2, 4.6, 1133)
  ,@(25, 21, 'Iron Skillet Restaurant', 2, 4.5, 470)
  ,@(26, 54, 'Kuma''s Corner', 2, 4.7, 2791)
  ,@(27, 16, 'Livery', 2, 4.7, 1490)
  ,@(28, 8, 'Maggiano''s Little Italy', 2, 4.4, 2257)
  ,@(29, 6, 'Major Restaurant', 2, 4.6, 427)
  ,@(30, 42, 'Mama Carolla''s', 2, 4.7, 1639)
  ,@(31, 9, 'Meridian Restaurant & Bar', 3, 4.5, 365)
  ,@(32, 40, 'Mimi Blue Restaurants', 2, 4.6, 1054)
  ,@(33, 38, 'Nada', 2, 4.4, 1952)
  ,@(34, 24, 'Nesso', $null, 4.7, 219)
  ,@(35, 35, 'Ocean Prime', 4, 4.6, 958)
  ,@(36, 13, 'Olive Garden Italian Restaurant', 2, 4.4, 1395)
  ,@(37, 59, 'Papa Fattoush Restaurant', 1, 4.5, 344)
  ,@(38, 32, 'Perkins Restaurant & Bakery', 2, 4.3, 981)
  ,@(39, 3, 'Rusty Bucket Restaurant and Tavern', 2, 4.4, 946)
  ,@(40, 20, 'Ruth''s Chris Steak House', 4, 4.6, 969)
  ,@(41, 17, 'Sahm''s Restaurant', 2, 4.5, 793)
  ,@(42, 2, 'Seasons 52', 2, 4.5, 1339)
  ,@(43, 12, 'Sero''s Family Restaurant', 2, 4.5, 1158)
  ,@(44, 51, 'Slapfish', 2, 4.6, 317)
  ,@(45, 1, 'The Capital Grille', 4, 4.6, 821)
  ,@(46, 33, 'The Cheesecake Factory', 2, 4.2, 3306)
  ,@(47, 46, 'The Italian House on Park', 2, 4.8, 544)
  ,@(48, 50, 'The Oceanaire Seafood Room', 3, 4.5, 975)
  ,@(49, 55, 'The Old Spaghetti Factory', 2, 4.4, 2772)
  ,@(50, 47, 'The Rathskeller', 2, 4.5, 2685)
  ,@(51, 15, 'Tinker Street Restaurant', 3, 4.7, 702)
  ,@(52, 45, 'Twin Peaks Restaurant', 2, 4.5, 3906)
  ,@(53, 30, 'Weber Grill Restaurant', 2, 4.2, 2307)
  ,@(54, 7, 'Yard House', 2, 4.4, 2326)
  ,@(55, 11, 'Yats', 1, 4.8, 1279)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

$ws.Range("A56").EntireRow.Delete()
